$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.0111856823266219
$ws.Range("C2").Value = 0.0134228187919463
$ws.Range("D2").Value = 0.938105891126025
$ws.Range("E2").Value = 0.0119313944817301
$ws.Range("F2").Value = 0.980611483967189
$ws.Range("G2").Value = 0.968680089485459
$ws.Range("H2").Value = 0.0111856823266219
$ws.Range("I2").Value = 0.00969425801640567
$ws.Range("J2").Value = 0.00149142431021626
$ws.Range("K2").Value = 0.991051454138702
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0.00447427293064877
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0.0380313199105145
$ws.Range("S2").Value = 0.029082774049217
$ws.Range("T2").Value = 0.00149142431021626
$ws.Range("U2").Value = 0.0134228187919463
$ws.Range("V2").Value = 0.0141685309470544
$ws.Range("W2").Value = 0.0104399701715138
$ws.Range("X2").Value = 0.00820283370618941
$ws.Range("B3").Value = 0.967188665175242
$ws.Range("C3").Value = 0.00671140939597315
$ws.Range("D3").Value = 0.0104399701715138
$ws.Range("E3").Value = 0.969425801640567
$ws.Range("F3").Value = 0.0126771066368382
$ws.Range("G3").Value = 0.00223713646532438
$ws.Range("H3").Value = 0.00149142431021626
$ws.Range("I3").Value = 0.000745712155108128
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0.00596569724086503
$ws.Range("L3").Value = 0.992542878448919
$ws.Range("M3").Value = 0.00372856077554064
$ws.Range("N3").Value = 0.00149142431021626
$ws.Range("O3").Value = 0.00149142431021626
$ws.Range("P3").Value = 0.997017151379567
$ws.Range("Q3").Value = 0.991051454138702
$ws.Range("R3").Value = 0.00596569724086503
$ws.Range("S3").Value = 0.954511558538404
$ws.Range("T3").Value = 0.0186428038777032
$ws.Range("U3").Value = 0.000745712155108128
$ws.Range("V3").Value = 0
$ws.Range("W3").Value = 0.000745712155108128
$ws.Range("X3").Value = 0.0052199850857569
$ws.Range("B4").Value = 0.00298284862043251
$ws.Range("C4").Value = 0.00894854586129754
$ws.Range("D4").Value = 0.046234153616704
$ws.Range("E4").Value = 0.0052199850857569
$ws.Range("F4").Value = 0.00149142431021626
$ws.Range("G4").Value = 0.0275913497390007
$ws.Range("H4").Value = 0.986577181208054
$ws.Range("I4").Value = 0.98806860551827
$ws.Range("J4").Value = 0.997762863534676
$ws.Range("K4").Value = 0.00298284862043251
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0.000745712155108128
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0.000745712155108128
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0.00596569724086503
$ws.Range("R4").Value = 0.95600298284862
$ws.Range("S4").Value = 0.00149142431021626
$ws.Range("T4").Value = 0.00223713646532438
$ws.Range("U4").Value = 0.982848620432513
$ws.Range("V4").Value = 0.979865771812081
$ws.Range("W4").Value = 0.987322893363162
$ws.Range("X4").Value = 0.985831469052946
$ws.Range("B5").Value = 0.0186428038777032
$ws.Range("C5").Value = 0.970917225950783
$ws.Range("D5").Value = 0.0052199850857569
$ws.Range("E5").Value = 0.0134228187919463
$ws.Range("F5").Value = 0.0052199850857569
$ws.Range("G5").Value = 0.00149142431021626
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0.00149142431021626
$ws.Range("J5").Value = 0.000745712155108128
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0.00671140939597315
$ws.Range("M5").Value = 0.991051454138702
$ws.Range("N5").Value = 0.998508575689784
$ws.Range("O5").Value = 0.997762863534676
$ws.Range("P5").Value = 0.00298284862043251
$ws.Range("Q5").Value = 0.00298284862043251
$ws.Range("R5").Value = 0
$ws.Range("S5").Value = 0.0149142431021626
$ws.Range("T5").Value = 0.977628635346756
$ws.Range("U5").Value = 0.00223713646532438
$ws.Range("V5").Value = 0.0052199850857569
$ws.Range("W5").Value = 0.000745712155108128
$ws.Range("X5").Value = 0.000745712155108128
